# Auto-generated edit script: apply scheduled-runner market-data updates
# to the Zalera_Profits workbook (columns H..N per leve row, per sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Range("H116").Value = 3011.875
$ws.Range("I116").Value = 2866.1667
$ws.Range("J116").Value = 3449
$ws.Range("K116").Value = 2866.1667
$ws.Range("L116").Value = 3449
$ws.Range("M116").Value = 575.8332999999998
$ws.Range("N116").Value = -10333

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 22729022
$ws.Range("I2").Value = 45455544
$ws.Range("K2").Value = 45455544
$ws.Range("M2").Value = -45455431

# Row 25
$ws.Range("H25").Value = 1095.8334
$ws.Range("I25").Value = 931.8
$ws.Range("K25").Value = 931.8
$ws.Range("M25").Value = -529.8

# Row 103
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344

# Row 116
$ws.Range("H116").Value = 22729022
$ws.Range("I116").Value = 45455544
$ws.Range("K116").Value = 45455544
$ws.Range("M116").Value = -45453250

# Row 122
$ws.Range("H122").Value = 1548.9
$ws.Range("I122").Value = 1232.5
$ws.Range("K122").Value = 3697.5
$ws.Range("M122").Value = -1247.5

# Row 132
$ws.Range("H132").Value = 4131.75
$ws.Range("I132").Value = 3086.4783
$ws.Range("K132").Value = 9259.4349
$ws.Range("M132").Value = -6729.4349

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 22729022
$ws.Range("I3").Value = 45455544
$ws.Range("K3").Value = 45455544
$ws.Range("M3").Value = -45455430

# Row 21
$ws.Range("H21").Value = 40385.5
$ws.Range("J21").Value = 40385.5
$ws.Range("L21").Value = 40385.5
$ws.Range("N21").Value = -40857.5

# Row 86
$ws.Range("H86").Value = 122068.53
$ws.Range("I86").Value = 4684.3335
$ws.Range("K86").Value = 4684.3335
$ws.Range("M86").Value = -3561.3335

# Row 89
$ws.Range("H89").Value = 122068.53
$ws.Range("I89").Value = 4684.3335
$ws.Range("K89").Value = 23421.6675
$ws.Range("M89").Value = -17805.6675

# Row 100
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164

# Row 103
$ws.Range("H103").Value = 200030100
$ws.Range("J103").Value = 200030100
$ws.Range("L103").Value = 200030100
$ws.Range("N103").Value = -200032444

# Row 133
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120

# Row 134
$ws.Range("H134").Value = 3181.5715
$ws.Range("I134").Value = 1498.6818
$ws.Range("J134").Value = 9352.166999999999
$ws.Range("K134").Value = 4496.0454
$ws.Range("L134").Value = 28056.501
$ws.Range("M134").Value = -1961.0454
$ws.Range("N134").Value = -33126.501

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 35718220
$ws.Range("I31").Value = 111113030
$ws.Range("J31").Value = 4885.8423
$ws.Range("K31").Value = 111113030
$ws.Range("L31").Value = 4885.8423
$ws.Range("M31").Value = -111112735
$ws.Range("N31").Value = -5475.8423

# Row 34
$ws.Range("H34").Value = 35718220
$ws.Range("I34").Value = 111113030
$ws.Range("J34").Value = 4885.8423
$ws.Range("K34").Value = 111113030
$ws.Range("L34").Value = 4885.8423
$ws.Range("M34").Value = -111112828
$ws.Range("N34").Value = -5289.8423

# Row 69
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 5000
$ws.Range("M69").Value = -4251

# Row 72
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 15000
$ws.Range("M72").Value = -11256

# Row 92
$ws.Range("H92").Value = 30999.5
$ws.Range("J92").Value = 30999.5
$ws.Range("L92").Value = 30999.5
$ws.Range("N92").Value = -35991.5

# Row 96
$ws.Range("H96").Value = 39660.145
$ws.Range("J96").Value = 39660.145
$ws.Range("L96").Value = 39660.145
$ws.Range("N96").Value = -45152.145

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 4189.6
$ws.Range("I55").Value = 316.66666
$ws.Range("K55").Value = 949.9999799999999
$ws.Range("M55").Value = -772.9999799999999

# Row 115
$ws.Range("H115").Value = 11165.5
$ws.Range("I115").Value = 2300
$ws.Range("K115").Value = 6900
$ws.Range("M115").Value = -5725

# Row 117
$ws.Range("H117").Value = 335599.8
$ws.Range("J117").Value = 335599.8
$ws.Range("L117").Value = 1006799.4
$ws.Range("N117").Value = -1013683.4

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 72297.71000000001
$ws.Range("I2").Value = 1004.1818
$ws.Range("J2").Value = 333707.34
$ws.Range("K2").Value = 1004.1818
$ws.Range("L2").Value = 333707.34
$ws.Range("M2").Value = -891.1818
$ws.Range("N2").Value = -333933.34

# Row 98
$ws.Range("H98").Value = 25998.5
$ws.Range("L98").Value = 25998.5
$ws.Range("N98").Value = -31988.5

# Row 132
$ws.Range("H132").Value = 4576.9355
$ws.Range("I132").Value = 2418.2
$ws.Range("K132").Value = 7254.599999999999
$ws.Range("M132").Value = -4724.599999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2211.3333
$ws.Range("I16").Value = 1353.4
$ws.Range("J16").Value = 6501
$ws.Range("K16").Value = 1353.4
$ws.Range("L16").Value = 6501
$ws.Range("M16").Value = -1183.4
$ws.Range("N16").Value = -6841

# Row 97
$ws.Range("H97").Value = 14886.833
$ws.Range("J97").Value = 14886.833
$ws.Range("L97").Value = 14886.833
$ws.Range("N97").Value = -16868.833

# Row 122
$ws.Range("H122").Value = 6786.091
$ws.Range("I122").Value = 4956.6665
$ws.Range("K122").Value = 14869.9995
$ws.Range("M122").Value = -12419.9995

# Row 132
$ws.Range("H132").Value = 6156.35
$ws.Range("I132").Value = 5233.846
$ws.Range("J132").Value = 7869.5713
$ws.Range("K132").Value = 15701.538
$ws.Range("L132").Value = 23608.7139
$ws.Range("M132").Value = -13171.538
$ws.Range("N132").Value = -28668.7139

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 6280.8
$ws.Range("J14").Value = 10100
$ws.Range("L14").Value = 10100
$ws.Range("N14").Value = -10436

# Row 20
$ws.Range("H20").Value = 72498.5
$ws.Range("J20").Value = 72498.5
$ws.Range("L20").Value = 72498.5
$ws.Range("N20").Value = -72978.5

# Row 132
$ws.Range("H132").Value = 4867.972
$ws.Range("I132").Value = 2373
$ws.Range("J132").Value = 7986.6875
$ws.Range("K132").Value = 7119
$ws.Range("L132").Value = 23960.0625
$ws.Range("M132").Value = -4589
$ws.Range("N132").Value = -29020.0625

# Row 133
$ws.Range("H133").Value = 69924.25
$ws.Range("J133").Value = 69924.25
$ws.Range("L133").Value = 69924.25
$ws.Range("N133").Value = -80044.25
